# Daily attendance processing - 2026-01-14 11:08:54
# Swap the order of "System" and the recorder's email address in the
# "Recorded By" column (G) wherever the combined "System, <email>" value
# is present, turning it into "<email>, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value2
    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
